# Auto-applies the BRVM automated update (GitHub Actions) to the
# "Recommandations" and "Top_YTD" sheets: cell values are refreshed
# in place and the now-unused trailing rows are cleared so the used
# range on "Recommandations" shrinks from A1:G45 to A1:G40.

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

$recoData = @{
    2 = @('BRVM - SERVICES PUBLICS', 0, 6, 2550.03, 115.08, '🟡 Observer', '➖ Neutre')
    3 = @('NEI-CEDA CI', 0, 3, 2200, 770, '🟡 Observer', '➖ Neutre')
    4 = @('AIR LIQUIDE CI', 0, 3, 2110, 700, '🟡 Observer', '➖ Neutre')
    5 = @('BRVM - DISTRIBUTION', 0, 3, 1824.47, 606.72, '🟡 Observer', '➖ Neutre')
    6 = @('BRVM - AUTRES SECTEURS', 0, 3, 1809.61, 605.27, '🟡 Observer', '➖ Neutre')
    7 = @('BRVM - TRANSPORT', 0, 3, 1079.45, 359.82, '🟡 Observer', '➖ Neutre')
    8 = @('BRVM - AGRICULTURE', 0, 3, 1063.34, 346.24, '🟡 Observer', '➖ Neutre')
    9 = @('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 3, 694.12, 229.76, '🟡 Observer', '➖ Neutre')
    10 = @('BRVM-PRESTIGE', 0, 3, 431.51, 143.66, '🟡 Observer', '➖ Neutre')
    11 = @('BRVM - FINANCES', 0, 3, 430.8, 143.59, '🟡 Observer', '➖ Neutre')
    12 = @('BRVM-PRINCIPAL                    (**)', 0, 2, 423.48, 213.74, '🟡 Observer', '➖ Neutre')
    13 = @('BRVM - SERVICES FINANCIERS', 0, 3, 423.39, 141.12, '🟡 Observer', '➖ Neutre')
    14 = @('BRVM - INDUSTRIELS', 0, 3, 376.72, 123.7, '🟡 Observer', '➖ Neutre')
    15 = @('BRVM - ENERGIE', 0, 3, 322.55, 108.19, '🟡 Observer', '➖ Neutre')
    16 = @('BRVM - TELECOMMUNICATIONS', 0, 3, 287.49, 95.46, '🟡 Observer', '➖ Neutre')
    17 = @('BRVM - INDUSTRIE                (**)', 0, 1, 235.98, 235.98, '🟡 Observer', '➖ Neutre')
    18 = @('BRVM - INDUSTRIE                  (**)', 0, 1, 218.47, 218.47, '🟡 Observer', '➖ Neutre')
    19 = @('BRVM - CONSOMMATION DE BASE         (**)', 0, 1, 202.47, 202.47, '🟡 Observer', '➖ Neutre')
    20 = @('BRVM - CONSOMMATION DE BASE          (**)', 0, 1, 193.64, 193.64, '🟡 Observer', '➖ Neutre')
    21 = @('NEI-CEDA CI (NEIC)', 3, 0, 16.94, 7.14, '🟢 Achat', '✅ Renforcer')
    22 = @('SOLIBRA CI (SLBC)', 2, 0, 14.98, 7.49, '🟡 Observer', '➖ Neutre')
    23 = @('CFAO MOTORS CI (CFAC)', 2, 0, 14.51, 7.08, '🟡 Observer', '➖ Neutre')
    24 = @('UNILEVER CI (UNLC)', 2, 0, 14.48, 6.98, '🟡 Observer', '➖ Neutre')
    25 = @('SICABLE CI (CABC)', 1, 0, 7.25, 7.25, '🟡 Observer', '➖ Neutre')
    26 = @('AFRICA GLOBAL LOGISTICS CI (SDSC)', 1, 0, 2.76, 2.76, '🟡 Observer', '➖ Neutre')
    27 = @('SOCIETE IVOIRIENNE DE BANQUE  (SIBC)', 1, 0, 2.59, 2.59, '🟡 Observer', '➖ Neutre')
    28 = @('TOTAL', 0, 3, 0, 0, '🟡 Observer', '➖ Neutre')
    29 = @('TRACTAFRIC MOTORS CI (PRSC)', 1, 1, -0.01, -7.5, '🟡 Observer', '👀 À surveiller')
    30 = @('SAFCA CI (SAFC)', 1, 1, -0.11, -7.5, '🟡 Observer', '👀 À surveiller')
    31 = @('TOTALENERGIES MARKETING CI (TTLC)', 0, 1, -1.67, -1.67, '🟡 Observer', '➖ Neutre')
    32 = @('SAPH CI (SPHC)', 1, 1, -2, 2.53, '🟡 Observer', '👀 À surveiller')
    33 = @('ORANGE COTE D''IVOIRE (ORAC)', 0, 1, -2.03, -2.03, '🟡 Observer', '➖ Neutre')
    34 = @('SOGB CI (SOGC)', 0, 1, -3.83, -3.83, '🟡 Observer', '➖ Neutre')
    35 = @('ECOBANK TRANS. INCORP. TG (ETIT)', 0, 1, -4.55, -4.55, '🟡 Observer', '➖ Neutre')
    36 = @('LOTERIE NATIONALE DU BENIN (LNBB)', 0, 2, -6.18, -4.63, '🟡 Observer', '➖ Neutre')
    37 = @('PALM CI (PALC)', 0, 1, -7.45, -7.45, '🟡 Observer', '➖ Neutre')
    38 = @('BANK OF AFRICA ML (BOAM)', 0, 1, -7.47, -7.47, '🟡 Observer', '➖ Neutre')
    39 = @('NESTLE CI (NTLC)', 0, 1, -7.48, -7.48, '🟡 Observer', '➖ Neutre')
    40 = @('FILTISAC CI (FTSC)', 0, 3, -17.76, -7.5, '🔴 Vente', '⚠️ Risque de décrochage')
}

foreach ($r in $recoData.Keys) {
    $row = $recoData[$r]
    $wsReco.Cells.Item($r, 1).Value = $row[0]
    $wsReco.Cells.Item($r, 2).Value = $row[1]
    $wsReco.Cells.Item($r, 3).Value = $row[2]
    $wsReco.Cells.Item($r, 4).Value = $row[3]
    $wsReco.Cells.Item($r, 5).Value = $row[4]
    $wsReco.Cells.Item($r, 6).Value = $row[5]
    $wsReco.Cells.Item($r, 7).Value = $row[6]
}

# Rows 41-45 no longer exist after the refresh; clear them so the
# sheet dimension shrinks back to A1:G40, matching the source feed.
$wsReco.Range("A41:G45").ClearContents()

$ytdData = @{
    2 = @('BRVM - SERVICES PUBLICS', 569575.63)
    3 = @('NEI-CEDA CI', 57652.77)
    4 = @('AIR LIQUIDE CI', 51740)
    5 = @('BRVM - DISTRIBUTION', 35369.04)
    6 = @('BRVM - AUTRES SECTEURS', 34672.76)
    7 = @('BRVM - TRANSPORT', 9620.86)
    8 = @('BRVM - AGRICULTURE', 9277.54)
    9 = @('BRVM - CONSOMMATION DISCRETIONNAIRE', 3533.79)
    10 = @('BRVM-PRESTIGE', 1349.74)
    11 = @('BRVM - FINANCES', 1345.54)
}

foreach ($r in $ytdData.Keys) {
    $row = $ytdData[$r]
    $wsYtd.Cells.Item($r, 1).Value = $row[0]
    $wsYtd.Cells.Item($r, 2).Value = $row[1]
}

